# Add J088 data from Yuanyuan to the "Tidy Data" sheet.
#
# Layout being built (rows 269-288 on "Tidy Data"):
#   269            : "DARK" / "DU mol%" / "DTh mol%" section header (like row 1, 17, 42...)
#   270-278 (9 rows): J088 dark-zone DU/DTh mol% measurements
#   279            : "light" sub-header (single cell, like rows 9, 97...)
#   280-288 (9 rows): J088 light-zone DU/DTh mol% measurements
#
# The old trailing "max"/"min" summary rows (old rows 270-271) are removed
# since the table no longer ends there.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Tidy Data")

# ------------------------------------------------------------------
# 1. Make room: insert a fresh row at 269 (old row 270 "max" shifts to
#    271, old row 271 "min" shifts to 272), then delete those two old
#    summary rows outright since they are being replaced by new data.
# ------------------------------------------------------------------
$ws2.Rows.Item(269).Insert()
$ws2.Range("A271:H272").Delete(-4162)

# ------------------------------------------------------------------
# 2. Row 269 header.
# ------------------------------------------------------------------
$ws2.Range("A269").Value = "DARK"
$ws2.Range("B269").Value = "DU mol%"
$ws2.Range("C269").Value = "DTh mol%"

# ------------------------------------------------------------------
# 3. DARK block: rows 270-278, sample J088.
# ------------------------------------------------------------------
$ws2.Range("A270").Value = "J088"
$ws2.Range("B270").Value = 7.6044379255531567
$ws2.Range("C270").Value = 4.9940033636871544
$ws2.Range("E270").Value = "J088"
$ws2.Range("F270").Value = "dark"

$ws2.Range("A271").Value = "J088"
$ws2.Range("B271").Value = 7.449855047858768
$ws2.Range("C271").Value = 4.982827504310813
$ws2.Range("E271").Value = "J088"
$ws2.Range("F271").Value = "dark"

$ws2.Range("A272").Value = "J088"
$ws2.Range("B272").Value = 4.9805936583139614
$ws2.Range("C272").Value = 3.0603781441021165
$ws2.Range("E272").Value = "J088"
$ws2.Range("F272").Value = "dark"

$ws2.Range("A273").Value = "J088"
$ws2.Range("B273").Value = 4.2841936663886031
$ws2.Range("C273").Value = 2.7233242815024612
$ws2.Range("E273").Value = "J088"
$ws2.Range("F273").Value = "dark"

$ws2.Range("A274").Value = "J088"
$ws2.Range("B274").Value = 4.916394520179284
$ws2.Range("C274").Value = 2.8869535542142217
$ws2.Range("E274").Value = "J088"
$ws2.Range("F274").Value = "dark"

$ws2.Range("A275").Value = "J088"
$ws2.Range("B275").Value = 4.4334132875956591
$ws2.Range("C275").Value = 2.7964247085539609
$ws2.Range("E275").Value = "J088"
$ws2.Range("F275").Value = "dark"

$ws2.Range("A276").Value = "J088"
$ws2.Range("B276").Value = 4.9307824446764785
$ws2.Range("C276").Value = 3.4224743164471358
$ws2.Range("E276").Value = "J088"
$ws2.Range("F276").Value = "dark"

$ws2.Range("A277").Value = "J088"
$ws2.Range("B277").Value = 4.7188352588500786
$ws2.Range("C277").Value = 3.241562268525974
$ws2.Range("E277").Value = "J088"
$ws2.Range("F277").Value = "dark"

$ws2.Range("A278").Value = "J088"
$ws2.Range("B278").Value = 4.7880528141735992
$ws2.Range("C278").Value = 3.0574297124086187
$ws2.Range("E278").Value = "J088"
$ws2.Range("F278").Value = "dark"

# ------------------------------------------------------------------
# 4. Row 279: "light" sub-header (single cell, column A only).
# ------------------------------------------------------------------
$ws2.Range("A279").Value = "light"

# ------------------------------------------------------------------
# 5. light block: rows 280-288, sample J088.
# ------------------------------------------------------------------
$ws2.Range("A280").Value = "J088"
$ws2.Range("B280").Value = 3.049237015043655
$ws2.Range("C280").Value = 1.5495900311668434
$ws2.Range("E280").Value = "J088"
$ws2.Range("F280").Value = "light"

$ws2.Range("A281").Value = "J088"
$ws2.Range("B281").Value = 2.9926175205895444
$ws2.Range("C281").Value = 1.3909776813414279
$ws2.Range("E281").Value = "J088"
$ws2.Range("F281").Value = "light"

$ws2.Range("A282").Value = "J088"
$ws2.Range("B282").Value = 3.5940546942409966
$ws2.Range("C282").Value = 1.6952622648727274
$ws2.Range("E282").Value = "J088"
$ws2.Range("F282").Value = "light"

$ws2.Range("A283").Value = "J088"
$ws2.Range("B283").Value = 4.4072763009402971
$ws2.Range("C283").Value = 2.2108643053011772
$ws2.Range("E283").Value = "J088"
$ws2.Range("F283").Value = "light"

$ws2.Range("A284").Value = "J088"
$ws2.Range("B284").Value = 0.3496203140788654
$ws2.Range("C284").Value = 0.17164864874609742
$ws2.Range("E284").Value = "J088"
$ws2.Range("F284").Value = "light"

$ws2.Range("A285").Value = "J088"
$ws2.Range("B285").Value = 0.38976938956873225
$ws2.Range("C285").Value = 0.19281644344871815
$ws2.Range("E285").Value = "J088"
$ws2.Range("F285").Value = "light"

$ws2.Range("A286").Value = "J088"
$ws2.Range("B286").Value = 2.9094257140474036
$ws2.Range("C286").Value = 1.4459182047490011
$ws2.Range("E286").Value = "J088"
$ws2.Range("F286").Value = "light"

$ws2.Range("A287").Value = "J088"
$ws2.Range("B287").Value = 2.5072417393647477
$ws2.Range("C287").Value = 1.0653441113007633
$ws2.Range("E287").Value = "J088"
$ws2.Range("F287").Value = "light"

$ws2.Range("A288").Value = "J088"
$ws2.Range("B288").Value = 2.6561184277004948
$ws2.Range("C288").Value = 0.98700524942835999
$ws2.Range("E288").Value = "J088"
$ws2.Range("F288").Value = "light"

# ------------------------------------------------------------------
# 6. Formulas in G/H (mirror B/C). First row of each block is a plain
#    formula; the rest of the block is filled as one range so the
#    writer groups them as shared formulas, exactly like the rest of
#    the sheet already does.
# ------------------------------------------------------------------
$ws2.Range("G270").Formula = "=B270"
$ws2.Range("H270").Formula = "=C270"
$ws2.Range("G271:G278").Formula = "=B271"
$ws2.Range("H271:H278").Formula = "=C271"

$ws2.Range("G280").Formula = "=B280"
$ws2.Range("H280").Formula = "=C280"
$ws2.Range("G281:G288").Formula = "=B281"
$ws2.Range("H281:H288").Formula = "=C281"

# ------------------------------------------------------------------
# 7. Styling to match the rest of the table:
#      - A269 uses the bold "DARK" section-header style (same as A1).
#      - E269 is an empty cell carrying the shaded-section style.
#      - A/E of every new data row (270-278, 280-288) carry the
#        shaded-section style used throughout the table (e.g. A7/E7).
# ------------------------------------------------------------------
$ws2.Range("A1").Copy() | Out-Null
$ws2.Range("A269").PasteSpecial(-4122) | Out-Null

$ws2.Range("A7").Copy() | Out-Null
$ws2.Range("E269").PasteSpecial(-4122) | Out-Null
$ws2.Range("A270:A278").PasteSpecial(-4122) | Out-Null
$ws2.Range("A280:A288").PasteSpecial(-4122) | Out-Null
$ws2.Range("E270:E278").PasteSpecial(-4122) | Out-Null
$ws2.Range("E280:E288").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 8. Sheet views: scroll / zoom / selection, set Sheet1 first and
#    "Tidy Data" last so "Tidy Data" ends up the active tab again.
# ------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.Zoom = 140
$ws1.Range("C280").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 160
$ws2.Range("A290").Select()
